$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 15
$ws.Range("B3").Value = "test2"
$ws.Range("C3").Value = "test2"
$ws.Range("D3").Value = "test2"
$ws.Range("E3").Value = "test2"
$ws.Range("F3").Value = "test2"
$ws.Range("G3").Value = 1

$ws.Range("G4").Select()
